$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 842.5
$ws.Range("I18").Value = 627.92
$ws.Range("J18").Value = 2630.6667
$ws.Range("K18").Value = 627.92
$ws.Range("L18").Value = 2630.6667
$ws.Range("M18").Value = -343.92
$ws.Range("N18").Value = -3198.6667

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 325634.66
$ws.Range("I33").Value = 457.8
$ws.Range("J33").Value = 867596.1
$ws.Range("K33").Value = 457.8
$ws.Range("L33").Value = 867596.1
$ws.Range("M33").Value = -228.8
$ws.Range("N33").Value = -868054.1

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 52780.65
$ws.Range("I64").Value = 127545
$ws.Range("J64").Value = 2937.75
$ws.Range("K64").Value = 127545
$ws.Range("L64").Value = 2937.75
$ws.Range("M64").Value = -127297
$ws.Range("N64").Value = -3433.75

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 52780.65
$ws.Range("I67").Value = 127545
$ws.Range("J67").Value = 2937.75
$ws.Range("K67").Value = 127545
$ws.Range("L67").Value = 2937.75
$ws.Range("M67").Value = -126687
$ws.Range("N67").Value = -4653.75

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4300.4
$ws.Range("I76").Value = 4000.8
$ws.Range("J76").Value = 4600
$ws.Range("K76").Value = 4000.8
$ws.Range("L76").Value = 4600
$ws.Range("M76").Value = -3685.8
$ws.Range("N76").Value = -5230

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4300.4
$ws.Range("I79").Value = 4000.8
$ws.Range("J79").Value = 4600
$ws.Range("K79").Value = 4000.8
$ws.Range("L79").Value = 4600
$ws.Range("M79").Value = -2908.8
$ws.Range("N79").Value = -6784

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3557.842
$ws.Range("I116").Value = 3868.75
$ws.Range("J116").Value = 1899.6666
$ws.Range("K116").Value = 3868.75
$ws.Range("L116").Value = 1899.6666
$ws.Range("M116").Value = -426.75
$ws.Range("N116").Value = -8783.6666

# ARM row 14
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2535.4546
$ws.Range("I14").Value = 1296.6666
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 1296.6666
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -1121.6666
$ws.Range("N14").Value = -3350

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 814.1795
$ws.Range("I74").Value = 742.3125
$ws.Range("J74").Value = 1142.7142
$ws.Range("K74").Value = 742.3125
$ws.Range("L74").Value = 1142.7142
$ws.Range("M74").Value = 131.6875
$ws.Range("N74").Value = -2890.7142

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 814.1795
$ws.Range("I77").Value = 742.3125
$ws.Range("J77").Value = 1142.7142
$ws.Range("K77").Value = 3711.5625
$ws.Range("L77").Value = 5713.571
$ws.Range("M77").Value = 656.4375
$ws.Range("N77").Value = -14449.571

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1509.7778
$ws.Range("I88").Value = 1548
$ws.Range("J88").Value = 1433.3334
$ws.Range("K88").Value = 1548
$ws.Range("L88").Value = 1433.3334
$ws.Range("M88").Value = -1142
$ws.Range("N88").Value = -2245.3334

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1509.7778
$ws.Range("I91").Value = 1548
$ws.Range("J91").Value = 1433.3334
$ws.Range("K91").Value = 1548
$ws.Range("L91").Value = 1433.3334
$ws.Range("M91").Value = -144
$ws.Range("N91").Value = -4241.3334

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 15171.955
$ws.Range("I132").Value = 19579.281
$ws.Range("J132").Value = 4323.154
$ws.Range("K132").Value = 58737.84299999999
$ws.Range("L132").Value = 12969.462
$ws.Range("M132").Value = -56207.84299999999
$ws.Range("N132").Value = -18029.462

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2671.4285
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2671.4285
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2671.4285
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3919.4285

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2671.4285
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2671.4285
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 13357.1425
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -19597.1425

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9986.416999999999
$ws.Range("I99").Value = 1648.6666
$ws.Range("J99").Value = 18324.166
$ws.Range("K99").Value = 1648.6666
$ws.Range("L99").Value = 18324.166
$ws.Range("M99").Value = -150.6666
$ws.Range("N99").Value = -21320.166

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9986.416999999999
$ws.Range("I126").Value = 1648.6666
$ws.Range("J126").Value = 18324.166
$ws.Range("K126").Value = 4945.9998
$ws.Range("L126").Value = 54972.49800000001
$ws.Range("M126").Value = -2475.9998
$ws.Range("N126").Value = -59912.49800000001

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2050.6924
$ws.Range("I81").Value = 1124.875
$ws.Range("J81").Value = 3532
$ws.Range("K81").Value = 3374.625
$ws.Range("L81").Value = 10596
$ws.Range("M81").Value = -2251.625
$ws.Range("N81").Value = -12842

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2050.6924
$ws.Range("I84").Value = 1124.875
$ws.Range("J84").Value = 3532
$ws.Range("K84").Value = 10123.875
$ws.Range("L84").Value = 31788
$ws.Range("M84").Value = -4507.875
$ws.Range("N84").Value = -43020

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 101356
$ws.Range("I70").Value = 158796.61
$ws.Range("J70").Value = 8015
$ws.Range("K70").Value = 158796.61
$ws.Range("L70").Value = 8015
$ws.Range("M70").Value = -158526.61
$ws.Range("N70").Value = -8555

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 101356
$ws.Range("I73").Value = 158796.61
$ws.Range("J73").Value = 8015
$ws.Range("K73").Value = 158796.61
$ws.Range("L73").Value = 8015
$ws.Range("M73").Value = -157860.61
$ws.Range("N73").Value = -9887

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 100104664
$ws.Range("I80").Value = 200207360
$ws.Range("J80").Value = 1966
$ws.Range("K80").Value = 200207360
$ws.Range("L80").Value = 1966
$ws.Range("M80").Value = -200206362
$ws.Range("N80").Value = -3962

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 100104664
$ws.Range("I83").Value = 200207360
$ws.Range("J83").Value = 1966
$ws.Range("K83").Value = 1001036800
$ws.Range("L83").Value = 9830
$ws.Range("M83").Value = -1001031808
$ws.Range("N83").Value = -19814

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 202733.5
$ws.Range("I102").Value = 1469.4286
$ws.Range("J102").Value = 672349.7
$ws.Range("K102").Value = 1469.4286
$ws.Range("L102").Value = 672349.7
$ws.Range("M102").Value = 152.5714
$ws.Range("N102").Value = -675593.7

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4566.6665
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 5600
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 5600
$ws.Range("M7").Value = -2388
$ws.Range("N7").Value = -5824

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 258922.3
$ws.Range("I55").Value = 517312.1
$ws.Range("J55").Value = 532.5
$ws.Range("K55").Value = 517312.1
$ws.Range("L55").Value = 532.5
$ws.Range("M55").Value = -517139.1
$ws.Range("N55").Value = -878.5

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4566.6665
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 5600
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 16800
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -21740

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6924.7856
$ws.Range("I54").Value = 7070
$ws.Range("J54").Value = 6913.615
$ws.Range("K54").Value = 7070
$ws.Range("L54").Value = 6913.615
$ws.Range("M54").Value = -6550
$ws.Range("N54").Value = -7953.615

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 167400
$ws.Range("I81").Value = 333803.66
$ws.Range("J81").Value = 111932.11
$ws.Range("K81").Value = 667607.3199999999
$ws.Range("L81").Value = 223864.22
$ws.Range("M81").Value = -666546.3199999999
$ws.Range("N81").Value = -225986.22

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 167400
$ws.Range("I84").Value = 333803.66
$ws.Range("J84").Value = 111932.11
$ws.Range("K84").Value = 3338036.6
$ws.Range("L84").Value = 1119321.1
$ws.Range("M84").Value = -3332732.6
$ws.Range("N84").Value = -1129929.1

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2470.1177
$ws.Range("I122").Value = 1732.4667
$ws.Range("J122").Value = 8002.5
$ws.Range("K122").Value = 5197.4001
$ws.Range("L122").Value = 24007.5
$ws.Range("M122").Value = -2747.4001
$ws.Range("N122").Value = -28907.5
